# Colorado overview workbook update — "Update factsheets with text edits
# from COMM": the "No. of 990 Filers w/ Gov Grants" count column (and the
# grand-total cells) switch from numeric storage to literal/text storage
# (numbers formatted with thousands separators where applicable), the
# Baca County placeholder row on the County sheet gets real percentage /
# currency text, and a new "Total" row is appended to the County sheet.

function Set-TextCell {
    param($ws, $row, $col, [string]$text)
    # Force literal text storage: a leading apostrophe stops Excel from
    # re-parsing number-shaped strings ("97", "2,634", "$0", "0.00%", ...)
    # back into numeric/currency/percentage values. Re-stamping the style
    # to "Normal" afterwards strips the quote-prefix style Excel applies,
    # so the cell ends up with no explicit style index — matching a plain
    # literal-text cell.
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "Overall": A2 2634 -> "2,634" (text)
# ----------------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")
Set-TextCell $wsOverall 2 1 "2,634"

# ----------------------------------------------------------------------
# Sheet "County": B2:B63 numeric -> text (same value); Baca County (row
# 64) placeholder zeros -> percentage / currency text; append new row 65
# "Total" with the state-wide totals.
# ----------------------------------------------------------------------
$wsCounty = $wb.Worksheets.Item("County")

$countyCounts = @(
    @(2, "97"), @(3, "22"), @(4, "172"), @(5, "14"), @(6, "4"),
    @(7, "232"), @(8, "25"), @(9, "22"), @(10, "1"), @(11, "8"),
    @(12, "3"), @(13, "3"), @(14, "2"), @(15, "5"), @(16, "16"),
    @(17, "619"), @(18, "2"), @(19, "68"), @(20, "51"), @(21, "238"),
    @(22, "2"), @(23, "12"), @(24, "53"), @(25, "1"), @(26, "16"),
    @(27, "34"), @(28, "3"), @(29, "4"), @(30, "2"), @(31, "174"),
    @(32, "2"), @(33, "6"), @(34, "54"), @(35, "8"), @(36, "136"),
    @(37, "13"), @(38, "6"), @(39, "11"), @(40, "55"), @(41, "2"),
    @(42, "5"), @(43, "23"), @(44, "31"), @(45, "8"), @(46, "13"),
    @(47, "9"), @(48, "8"), @(49, "2"), @(50, "44"), @(51, "10"),
    @(52, "65"), @(53, "7"), @(54, "9"), @(55, "47"), @(56, "7"),
    @(57, "3"), @(58, "28"), @(59, "43"), @(60, "15"), @(61, "1"),
    @(62, "48"), @(63, "10")
)
foreach ($pair in $countyCounts) {
    Set-TextCell $wsCounty $pair[0] 2 $pair[1]
}

# Baca County (row 64): placeholder 0 values become formatted text.
Set-TextCell $wsCounty 64 2 "0.00%"
Set-TextCell $wsCounty 64 3 "`$0"
Set-TextCell $wsCounty 64 4 "0.00%"
Set-TextCell $wsCounty 64 5 "0.00%"
Set-TextCell $wsCounty 64 6 "0.00%"

# New row 65: state-wide Total, matching the other sheets' Total rows.
Set-TextCell $wsCounty 65 1 "Total"
Set-TextCell $wsCounty 65 2 "2,634"
Set-TextCell $wsCounty 65 3 "`$3,482,686,129"
Set-TextCell $wsCounty 65 4 "11.50%"
Set-TextCell $wsCounty 65 5 "-6.54%"
Set-TextCell $wsCounty 65 6 "59.83%"

# ----------------------------------------------------------------------
# Sheet "Congressional District": B2:B9 numeric -> text; B10 (Total)
# 2634 -> "2,634" (text).
# ----------------------------------------------------------------------
$wsCD = $wb.Worksheets.Item("Congressional District")

$cdCounts = @(
    @(2, "622"), @(3, "504"), @(4, "533"), @(5, "162"),
    @(6, "238"), @(7, "198"), @(8, "262"), @(9, "115")
)
foreach ($pair in $cdCounts) {
    Set-TextCell $wsCD $pair[0] 2 $pair[1]
}
Set-TextCell $wsCD 10 2 "2,634"

# ----------------------------------------------------------------------
# Sheet "Size": B2:B7 numeric -> text; B8 (Total) 2634 -> "2,634" (text).
# ----------------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")

$sizeCounts = @(
    @(2, "822"), @(3, "771"), @(4, "514"), @(5, "168"),
    @(6, "255"), @(7, "104")
)
foreach ($pair in $sizeCounts) {
    Set-TextCell $wsSize $pair[0] 2 $pair[1]
}
Set-TextCell $wsSize 8 2 "2,634"

# ----------------------------------------------------------------------
# Sheet "Subsector": B2:B13 numeric -> text; B14 (Total) 2634 -> "2,634"
# (text).
# ----------------------------------------------------------------------
$wsSubsector = $wb.Worksheets.Item("Subsector")

$subsectorCounts = @(
    @(2, "249"), @(3, "391"), @(4, "176"), @(5, "226"), @(6, "9"),
    @(7, "805"), @(8, "36"), @(9, "1"), @(10, "219"), @(11, "63"),
    @(12, "451"), @(13, "8")
)
foreach ($pair in $subsectorCounts) {
    Set-TextCell $wsSubsector $pair[0] 2 $pair[1]
}
Set-TextCell $wsSubsector 14 2 "2,634"
